# Apply the "cement data" update:
#  - Rename the single sheet "Sheet1" -> "data"
#  - Add a new sheet "cements" right after "data", populate it with the
#    Sample / Cement lookup table, and make it the active/selected sheet
#  - Restore the selection on "data" to match the post-edit state

$wb = $excel.ActiveWorkbook

# --- rename existing sheet -------------------------------------------------
$dataSheet = $wb.Worksheets.Item(1)
$dataSheet.Name = "data"

# --- add the new "cements" sheet right after "data" -------------------------
$cementSheet = $wb.Worksheets.Add($null, $dataSheet)
$cementSheet.Name = "cements"

# --- populate header + 27 sample rows ---------------------------------------
$sampleCement = @(
    @("Sample", "Cement"),
    @("FUWM 1 1", "Micrite (SSF)"),
    @("FUWM 1 2", "Micrite"),
    @("FUWM 1 3", "Micrite"),
    @("FUWM 1 5", "Laminar Calcrete"),
    @("FUWM 1 6", "Laminar Microbial"),
    @("FUWM 3.5 1", "Micrite (SSF)"),
    @("FUWM 3.5 2", "Laminar Calcrete"),
    @("FUWM 3.5 3", "Laminar Microbial"),
    @("FUWM 3.5 5", "Laminar Calcrete"),
    @("FUWM 3.5 H 1", "Laminar Microbial"),
    @("FUWM 3.5 H 2", "Laminar Microbial"),
    @("FUWM 3.5 H 3", "Micrite (SSF)"),
    @("FUWM 8.5 1", "Laminar Microbial"),
    @("FUWM 8.5 2", "Manganese"),
    @("FUWM 8.5 3", "Manganese"),
    @("FUWM 8.5 4", "Microbial (SSF)"),
    @("FUWM 16 B 2", "Micrite (SSF)"),
    @("FUWM 16 H 2", "Micrite (SSF)"),
    @("FUWM 16 H 3", "Microbial"),
    @("FUWM 16 H 4", "Micrite"),
    @("FUWM 16 T 4", "Microbial"),
    @("FUWM 16 T 5", "Laminar Microbial"),
    @("FUWM 16 T 6", "Microbial"),
    @("FUWM 16 T 7", "Micrite (SSF)"),
    @("FUWM 16 T 8", "Microbial"),
    @("FUWM 16 T 9", "Micrite")
)

for ($i = 0; $i -lt $sampleCement.Length; $i++) {
    $rowNum = $i + 1
    $pair = $sampleCement[$i]
    $cementSheet.Cells.Item($rowNum, 1).Value = $pair[0]
    $cementSheet.Cells.Item($rowNum, 2).Value = $pair[1]
}

# bestFit-ish width for column A (sample names)
$cementSheet.Columns.Item(1).ColumnWidth = 12.14

# --- view/selection bookkeeping ---------------------------------------------
# "data" keeps a plain selection (no frozen/topLeft scroll) at B32
[void]$dataSheet.Range("B32").Select()

# "cements" becomes the active sheet/tab with its own selection at D21
[void]$cementSheet.Activate()
[void]$cementSheet.Range("D21").Select()
